$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-28 Sunday", "2024-01-29 Monday"),
    @("468÷4=", "509÷4="),
    @("382÷2=", "961÷9="),
    @("753÷5=", "364÷9="),
    @("901÷8=", "693÷8="),
    @("816÷9=", "725÷2="),
    @("643÷3=", "766÷2="),
    @("160÷9=", "191÷9="),
    @("590÷3=", "998÷4="),
    @("540÷8=", "783÷7="),
    @("974÷4=", "548÷2="),
    @("125÷6=", "728÷4="),
    @("184÷8=", "818÷4="),
    @("748÷2=", "230÷6="),
    @("146÷2=", "290÷8="),
    @("386÷5=", "176÷9="),
    @("920÷9=", "203÷2="),
    @("675÷8=", "563÷5="),
    @("804÷5=", "453÷9="),
    @("673÷9=", "265÷9="),
    @("869÷4=", "948÷4="),
    @("458÷6=", "362÷9="),
    @("100÷5=", "724÷5="),
    @("971÷5=", "828÷2="),
    @("385÷6=", "887÷7="),
    @("229÷7=", "678÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
